# Add two new columns (I = "I0", J = "IF") to the worksheet, mirroring the
# styling already used by the existing header/data columns, and populate
# them with their per-row values for rows 2..66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy style from the existing header cell H1 so the new
# headers look consistent with the rest of the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row data for columns I ("I0") and J ("IF"), keyed by row number (2..66).
$data = @{
    2  = @(3, 3)
    3  = @(8, 9)
    4  = @(9, 9)
    5  = @(6, 6)
    6  = @(6, 6)
    7  = @(9, 9)
    8  = @(6, 7)
    9  = @(4, 4)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(6, 6)
    17 = @(6, 6)
    18 = @(5, 5)
    19 = @(5, 5)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(7, 7)
    24 = @(9, 9)
    25 = @(7, 8)
    26 = @(7, 7)
    27 = @(6, 7)
    28 = @(6, 6)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(6, 6)
    32 = @(6, 6)
    33 = @(3, 4)
    34 = @(6, 6)
    35 = @(6, 6)
    36 = @(6, 7)
    37 = @(4, 4)
    38 = @(6, 7)
    39 = @(6, 6)
    40 = @(5, 5)
    41 = @(5, 6)
    42 = @(6, 6)
    43 = @(7, 7)
    44 = @(7, 7)
    45 = @(3, 4)
    46 = @(6, 7)
    47 = @(10, 10)
    48 = @(6, 6)
    49 = @(5, 5)
    50 = @(6, 6)
    51 = @(8, 8)
    52 = @(7, 7)
    53 = @(7, 7)
    54 = @(7, 8)
    55 = @(6, 6)
    56 = @(2, 3)
    57 = @(8, 8)
    58 = @(7, 8)
    59 = @(9, 9)
    60 = @(5, 6)
    61 = @(6, 7)
    62 = @(9, 9)
    63 = @(7, 7)
    64 = @(5, 5)
    65 = @(3, 3)
    66 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($r, 10).Value = $vals[1]  # column J
}

$wb.Save()
